$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.346035666666667
$ws.Range("H2").Value = 7.038107
$ws.Range("I2").Value = 0.2908248950424448
$ws.Range("J2").Value = 0.2908248950424447
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.38723566666667
$ws.Range("N2").Value = 94.16170700000001
$ws.Range("O2").Value = 0.5539598599114094
$ws.Range("P2").Value = 0.5539598599114095
$ws.Range("Q2").Value = 73.63557435207213
$ws.Range("R2").Value = 662.720169168649
$ws.Range("S2").Value = 0.1611053181164631
$ws.Range("T2").Value = 0.1611053181164631

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.346035666666667
$ws.Range("H3").Value = 7.038107
$ws.Range("I3").Value = 0.2908248950424448
$ws.Range("J3").Value = 0.2908248950424447
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 17.55525033333333
$ws.Range("N3").Value = 52.665751
$ws.Range("O3").Value = 0.3098362697066353
$ws.Range("P3").Value = 0.3098362697066353
$ws.Range("Q3").Value = 41.18524341926189
$ws.Range("R3").Value = 370.667190773357
$ws.Range("S3").Value = 0.09010810061777483
$ws.Range("T3").Value = 0.09010810061777481

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.346035666666667
$ws.Range("H4").Value = 7.038107
$ws.Range("I4").Value = 0.2908248950424448
$ws.Range("J4").Value = 0.2908248950424447
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 7.717279333333334
$ws.Range("N4").Value = 23.151838
$ws.Range("O4").Value = 0.1362038703819552
$ws.Range("P4").Value = 0.1362038703819552
$ws.Range("Q4").Value = 18.10501256562956
$ws.Range("R4").Value = 162.945113090666
$ws.Range("S4").Value = 0.03961147630820688
$ws.Range("T4").Value = 0.03961147630820688

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.935458666666667
$ws.Range("H5").Value = 8.806376
$ws.Range("I5").Value = 0.3638923613841484
$ws.Range("J5").Value = 0.3638923613841484
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 31.38723566666667
$ws.Range("N5").Value = 94.16170700000001
$ws.Range("O5").Value = 0.5539598599114094
$ws.Range("P5").Value = 0.5539598599114095
$ws.Range("Q5").Value = 92.13593296042579
$ws.Range("R5").Value = 829.2233966438321
$ws.Range("S5").Value = 0.2015817615351948
$ws.Range("T5").Value = 0.2015817615351949

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.935458666666667
$ws.Range("H6").Value = 8.806376
$ws.Range("I6").Value = 0.3638923613841484
$ws.Range("J6").Value = 0.3638923613841484
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 17.55525033333333
$ws.Range("N6").Value = 52.665751
$ws.Range("O6").Value = 0.3098362697066353
$ws.Range("P6").Value = 0.3098362697066353
$ws.Range("Q6").Value = 51.53271173648623
$ws.Range("R6").Value = 463.794405628376
$ws.Range("S6").Value = 0.1127470518260034
$ws.Range("T6").Value = 0.1127470518260034

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.935458666666667
$ws.Range("H7").Value = 8.806376
$ws.Range("I7").Value = 0.3638923613841484
$ws.Range("J7").Value = 0.3638923613841484
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.717279333333334
$ws.Range("N7").Value = 23.151838
$ws.Range("O7").Value = 0.1362038703819552
$ws.Range("P7").Value = 0.1362038703819552
$ws.Range("Q7").Value = 22.65375450212089
$ws.Range("R7").Value = 203.883790519088
$ws.Range("S7").Value = 0.04956354802295016
$ws.Range("T7").Value = 0.04956354802295016

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.2753
$ws.Range("H8").Value = 3.8259
$ws.Range("I8").Value = 0.1580917945610786
$ws.Range("J8").Value = 0.1580917945610786
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 31.38723566666667
$ws.Range("N8").Value = 94.16170700000001
$ws.Range("O8").Value = 0.5539598599114094
$ws.Range("P8").Value = 0.5539598599114095
$ws.Range("Q8").Value = 40.0281416457
$ws.Range("R8").Value = 360.2532748113
$ws.Range("S8").Value = 0.08757650836819841
$ws.Range("T8").Value = 0.08757650836819846

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.2753
$ws.Range("H9").Value = 3.8259
$ws.Range("I9").Value = 0.1580917945610786
$ws.Range("J9").Value = 0.1580917945610786
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.55525033333333
$ws.Range("N9").Value = 52.665751
$ws.Range("O9").Value = 0.3098362697066353
$ws.Range("P9").Value = 0.3098362697066353
$ws.Range("Q9").Value = 22.3882107501
$ws.Range("R9").Value = 201.4938967509
$ws.Range("S9").Value = 0.04898257189803233
$ws.Range("T9").Value = 0.04898257189803234

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.2753
$ws.Range("H10").Value = 3.8259
$ws.Range("I10").Value = 0.1580917945610786
$ws.Range("J10").Value = 0.1580917945610786
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 7.717279333333334
$ws.Range("N10").Value = 23.151838
$ws.Range("O10").Value = 0.1362038703819552
$ws.Range("P10").Value = 0.1362038703819552
$ws.Range("Q10").Value = 9.8418463338
$ws.Range("R10").Value = 88.5766170042
$ws.Range("S10").Value = 0.02153271429484785
$ws.Range("T10").Value = 0.02153271429484785

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.839594
$ws.Range("H11").Value = 2.518782
$ws.Range("I11").Value = 0.1040797633205632
$ws.Range("J11").Value = 0.1040797633205632
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 31.38723566666667
$ws.Range("N11").Value = 94.16170700000001
$ws.Range("O11").Value = 0.5539598599114094
$ws.Range("P11").Value = 0.5539598599114095
$ws.Range("Q11").Value = 26.35253474231934
$ws.Range("R11").Value = 237.172812680874
$ws.Range("S11").Value = 0.05765601110867183
$ws.Range("T11").Value = 0.05765601110867184

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.839594
$ws.Range("H12").Value = 2.518782
$ws.Range("I12").Value = 0.1040797633205632
$ws.Range("J12").Value = 0.1040797633205632
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 17.55525033333333
$ws.Range("N12").Value = 52.665751
$ws.Range("O12").Value = 0.3098362697066353
$ws.Range("P12").Value = 0.3098362697066353
$ws.Range("Q12").Value = 14.73928284836467
$ws.Range("R12").Value = 132.653545635282
$ws.Range("S12").Value = 0.03224768561919279
$ws.Range("T12").Value = 0.03224768561919279

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.839594
$ws.Range("H13").Value = 2.518782
$ws.Range("I13").Value = 0.1040797633205632
$ws.Range("J13").Value = 0.1040797633205632
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.717279333333334
$ws.Range("N13").Value = 23.151838
$ws.Range("O13").Value = 0.1362038703819552
$ws.Range("P13").Value = 0.1362038703819552
$ws.Range("Q13").Value = 6.479381424590667
$ws.Range("R13").Value = 58.314432821316
$ws.Range("S13").Value = 0.01417606659269857
$ws.Range("T13").Value = 0.01417606659269857

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.6704439999999999
$ws.Range("H14").Value = 2.011332
$ws.Range("I14").Value = 0.08311118569176491
$ws.Range("J14").Value = 0.08311118569176491
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.38723566666667
$ws.Range("N14").Value = 94.16170700000001
$ws.Range("O14").Value = 0.5539598599114094
$ws.Range("P14").Value = 0.5539598599114095
$ws.Range("Q14").Value = 21.04338382930267
$ws.Range("R14").Value = 189.390454463724
$ws.Range("S14").Value = 0.04604026078288122
$ws.Range("T14").Value = 0.04604026078288123

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.6704439999999999
$ws.Range("H15").Value = 2.011332
$ws.Range("I15").Value = 0.08311118569176491
$ws.Range("J15").Value = 0.08311118569176491
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 17.55525033333333
$ws.Range("N15").Value = 52.665751
$ws.Range("O15").Value = 0.3098362697066353
$ws.Range("P15").Value = 0.3098362697066353
$ws.Range("Q15").Value = 11.76981225448133
$ws.Range("R15").Value = 105.928310290332
$ws.Range("S15").Value = 0.02575085974563192
$ws.Range("T15").Value = 0.02575085974563192

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.6704439999999999
$ws.Range("H16").Value = 2.011332
$ws.Range("I16").Value = 0.08311118569176491
$ws.Range("J16").Value = 0.08311118569176491
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.717279333333334
$ws.Range("N16").Value = 23.151838
$ws.Range("O16").Value = 0.1362038703819552
$ws.Range("P16").Value = 0.1362038703819552
$ws.Range("Q16").Value = 5.174003625357333
$ws.Range("R16").Value = 46.566032628216
$ws.Range("S16").Value = 0.01132006516325176
$ws.Range("T16").Value = 0.01132006516325176
